# Update the Lebanon MSME summary figures with higher-precision values.
# "Enterprises density (per 1000 people)" row (row 11): Micro / SMEs / MSMEs
# "Enterprises (% of total)" row (row 13): Micro / SMEs / MSMEs
#
# The source values are stored as text (not numbers) in the workbook, so we
# preserve that by quote-prefixing the new value before assigning it, and we
# restore each cell's original style afterwards so Excel's automatic
# "Number Stored as Text" formatting (quotePrefix) doesn't change the cell's
# visual style from what it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B11" = "41.51"
    "C11" = "2.73"
    "D11" = "44.25"
    "B13" = "93.53"
    "C13" = "6.16"
    "D13" = "99.68"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = $originalStyle
}
